$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 2 (shifts the existing data rows down by one) and strip
# any formatting it inherited from the header row above.
$ws.Rows(2).Insert()
$ws.Range("A2:F2").ClearFormats()

# New "registration" record values.
$ws.Range("A2").Value = 815
$ws.Range("B2").Value = "John"
$ws.Range("C2").Value = "Doe"
$ws.Range("D2").Value = "männlich"
$ws.Range("E2").Value = "Beerdigungen"

# Formatting that mirrors a small data-entry form: number right aligned,
# text fields left/top aligned, and a placeholder F column centered/top.
$ws.Range("A2").Font.Name = "Calibri"
$ws.Range("A2").HorizontalAlignment = -4152

$ws.Range("B2:E2").Font.Name = "Calibri"
$ws.Range("B2:E2").HorizontalAlignment = -4131
$ws.Range("B2:E2").VerticalAlignment = -4160

$ws.Range("F2").Font.Name = "Calibri"
$ws.Range("F2").HorizontalAlignment = -4108
$ws.Range("F2").VerticalAlignment = -4160

# New selection left behind by the edit (mirrors selecting the whole new row).
$ws.Rows(2).Select()
